$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 247.07143
$ws.Range("I2").Value = 174
$ws.Range("J2").Value = 344.5
$ws.Range("K2").Value = 174
$ws.Range("L2").Value = 344.5
$ws.Range("M2").Value = -61
$ws.Range("N2").Value = -570.5
$ws.Range("H33").Value = 239.16667
$ws.Range("I33").Value = 167.25
$ws.Range("K33").Value = 167.25
$ws.Range("M33").Value = 61.75
$ws.Range("H43").Value = 2184.5
$ws.Range("J43").Value = 2479.6667
$ws.Range("L43").Value = 2479.6667
$ws.Range("N43").Value = -2617.6667
$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3238
$ws.Range("H51").Value = 12198.143
$ws.Range("I51").Value = 40000
$ws.Range("K51").Value = 40000
$ws.Range("M51").Value = -39516
$ws.Range("H58").Value = 2612.375
$ws.Range("J58").Value = 20000
$ws.Range("L58").Value = 60000
$ws.Range("N58").Value = -60300
$ws.Range("H60").Value = 1000
$ws.Range("J60").Value = 1000
$ws.Range("L60").Value = 3000
$ws.Range("N60").Value = -3968
$ws.Range("H62").Value = 2882.48
$ws.Range("I62").Value = 3035.0625
$ws.Range("K62").Value = 3035.0625
$ws.Range("M62").Value = -2411.0625
$ws.Range("H65").Value = 2882.48
$ws.Range("I65").Value = 3035.0625
$ws.Range("K65").Value = 15175.3125
$ws.Range("M65").Value = -12055.3125
$ws.Range("H80").Value = 299.66666
$ws.Range("I80").Value = 285.625
$ws.Range("K80").Value = 856.875
$ws.Range("M80").Value = 141.125
$ws.Range("H82").Value = 8790.75
$ws.Range("J82").Value = 19909
$ws.Range("L82").Value = 59727
$ws.Range("N82").Value = -60539
$ws.Range("H83").Value = 299.66666
$ws.Range("I83").Value = 285.625
$ws.Range("K83").Value = 2570.625
$ws.Range("M83").Value = 2421.375
$ws.Range("H85").Value = 8790.75
$ws.Range("J85").Value = 19909
$ws.Range("L85").Value = 59727
$ws.Range("N85").Value = -62535
$ws.Range("H86").Value = 2713.8333
$ws.Range("I86").Value = 3175.8
$ws.Range("J86").Value = 404
$ws.Range("K86").Value = 3175.8
$ws.Range("L86").Value = 404
$ws.Range("M86").Value = -2052.8
$ws.Range("N86").Value = -2650
$ws.Range("H88").Value = 1112.091
$ws.Range("J88").Value = 1248.3334
$ws.Range("L88").Value = 1248.3334
$ws.Range("N88").Value = -2060.3334
$ws.Range("H89").Value = 2713.8333
$ws.Range("I89").Value = 3175.8
$ws.Range("J89").Value = 404
$ws.Range("K89").Value = 15879
$ws.Range("L89").Value = 2020
$ws.Range("M89").Value = -10263
$ws.Range("N89").Value = -13252
$ws.Range("H91").Value = 1112.091
$ws.Range("J91").Value = 1248.3334
$ws.Range("L91").Value = 1248.3334
$ws.Range("N91").Value = -4056.3334
$ws.Range("H92").Value = 1987.3334
$ws.Range("I92").Value = 1835.3
$ws.Range("K92").Value = 1835.3
$ws.Range("M92").Value = -587.3
$ws.Range("H98").Value = 5685.6665
$ws.Range("I98").Value = 1773.3
$ws.Range("K98").Value = 1773.3
$ws.Range("M98").Value = -275.3
$ws.Range("H106").Value = 5638.6665
$ws.Range("I106").Value = 6566.4
$ws.Range("K106").Value = 6566.4
$ws.Range("M106").Value = -5935.4
$ws.Range("H113").Value = 5214.4185
$ws.Range("I113").Value = 4001.12
$ws.Range("J113").Value = 6899.5557
$ws.Range("K113").Value = 4001.12
$ws.Range("L113").Value = 6899.5557
$ws.Range("M113").Value = -747.1199999999999
$ws.Range("N113").Value = -13407.5557
$ws.Range("H122").Value = 5685.6665
$ws.Range("I122").Value = 1773.3
$ws.Range("K122").Value = 5319.9
$ws.Range("M122").Value = -2869.9
$ws.Range("H125").Value = 2645.25
$ws.Range("I125").Value = 1981
$ws.Range("K125").Value = 17829
$ws.Range("M125").Value = -15369
$ws.Range("H132").Value = 29001.225
$ws.Range("I132").Value = 43573.52
$ws.Range("J132").Value = 4714.067
$ws.Range("K132").Value = 130720.56
$ws.Range("L132").Value = 14142.201
$ws.Range("M132").Value = -128190.56
$ws.Range("N132").Value = -19202.201
$ws.Range("H134").Value = 32912.082
$ws.Range("J134").Value = 32912.082
$ws.Range("L134").Value = 32912.082
$ws.Range("N134").Value = -43052.082
$ws.Range("H137").Value = 1792.1351
$ws.Range("I137").Value = 1277.12
$ws.Range("J137").Value = 2865.0833
$ws.Range("K137").Value = 3831.36
$ws.Range("L137").Value = 8595.249899999999
$ws.Range("M137").Value = -1281.36
$ws.Range("N137").Value = -13695.2499
$ws.Range("H138").Value = 6919.3335
$ws.Range("I138").Value = 6431.5835
$ws.Range("J138").Value = 7569.6665
$ws.Range("K138").Value = 19294.7505
$ws.Range("L138").Value = 22708.9995
$ws.Range("M138").Value = -14154.7505
$ws.Range("N138").Value = -32988.99950000001
$ws.Range("H141").Value = 3399.4443
$ws.Range("I141").Value = 2582.2068
$ws.Range("J141").Value = 6785.143
$ws.Range("K141").Value = 7746.6204
$ws.Range("L141").Value = 20355.429
$ws.Range("M141").Value = -2566.6204
$ws.Range("N141").Value = -30715.429

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1148.2
$ws.Range("I2").Value = 780.3333
$ws.Range("K2").Value = 780.3333
$ws.Range("M2").Value = -667.3333
$ws.Range("H25").Value = 2500
$ws.Range("I25").Value = 2500
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2500
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -2098
$ws.Range("N25").ClearContents()
$ws.Range("H32").Value = 2829.745
$ws.Range("I32").Value = 2690.8086
$ws.Range("J32").Value = 4462.25
$ws.Range("K32").Value = 2690.8086
$ws.Range("L32").Value = 4462.25
$ws.Range("M32").Value = -2403.8086
$ws.Range("N32").Value = -5036.25
$ws.Range("H61").Value = 7277.8125
$ws.Range("I61").Value = 1121.4166
$ws.Range("K61").Value = 1121.4166
$ws.Range("M61").Value = -909.4166
$ws.Range("H74").Value = 1314.6774
$ws.Range("I74").Value = 1139.7727
$ws.Range("K74").Value = 1139.7727
$ws.Range("M74").Value = -265.7727
$ws.Range("H77").Value = 1314.6774
$ws.Range("I77").Value = 1139.7727
$ws.Range("K77").Value = 5698.863499999999
$ws.Range("M77").Value = -1330.863499999999
$ws.Range("H88").Value = 11229.8
$ws.Range("I88").Value = 1300
$ws.Range("K88").Value = 1300
$ws.Range("M88").Value = -894
$ws.Range("H91").Value = 11229.8
$ws.Range("I91").Value = 1300
$ws.Range("K91").Value = 1300
$ws.Range("M91").Value = 104
$ws.Range("H92").Value = 55000
$ws.Range("J92").Value = 55000
$ws.Range("L92").Value = 55000
$ws.Range("N92").Value = -59992
$ws.Range("H102").Value = 1076.5333
$ws.Range("I102").Value = 1076.5333
$ws.Range("K102").Value = 1076.5333
$ws.Range("M102").Value = 545.4667
$ws.Range("H116").Value = 1148.2
$ws.Range("I116").Value = 780.3333
$ws.Range("K116").Value = 780.3333
$ws.Range("M116").Value = 1513.6667
$ws.Range("H132").Value = 2432.4688
$ws.Range("I132").Value = 1777.7778
$ws.Range("J132").Value = 5967.8
$ws.Range("K132").Value = 5333.3334
$ws.Range("L132").Value = 17903.4
$ws.Range("M132").Value = -2803.3334
$ws.Range("N132").Value = -22963.4
$ws.Range("H135").Value = 76949.25
$ws.Range("J135").Value = 76949.25
$ws.Range("L135").Value = 76949.25
$ws.Range("N135").Value = -87089.25
$ws.Range("H136").Value = 7277.8125
$ws.Range("I136").Value = 1121.4166
$ws.Range("K136").Value = 3364.2498
$ws.Range("M136").Value = -814.2498

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1148.2
$ws.Range("I3").Value = 780.3333
$ws.Range("K3").Value = 780.3333
$ws.Range("M3").Value = -666.3333
$ws.Range("H21").Value = 18699
$ws.Range("J21").Value = 18699
$ws.Range("L21").Value = 18699
$ws.Range("N21").Value = -19171
$ws.Range("H25").Value = 1128.5
$ws.Range("I25").Value = 1004.6667
$ws.Range("J25").Value = 1500
$ws.Range("K25").Value = 1004.6667
$ws.Range("L25").Value = 1500
$ws.Range("M25").Value = -769.6667
$ws.Range("N25").Value = -1970
$ws.Range("H64").Value = 755.2143
$ws.Range("J64").Value = 979.5714
$ws.Range("L64").Value = 979.5714
$ws.Range("N64").Value = -1429.5714
$ws.Range("H67").Value = 755.2143
$ws.Range("J67").Value = 979.5714
$ws.Range("L67").Value = 979.5714
$ws.Range("N67").Value = -2539.5714
$ws.Range("H86").Value = 37034.92
$ws.Range("I86").Value = 9556.75
$ws.Range("K86").Value = 9556.75
$ws.Range("M86").Value = -8433.75
$ws.Range("H89").Value = 37034.92
$ws.Range("I89").Value = 9556.75
$ws.Range("K89").Value = 47783.75
$ws.Range("M89").Value = -42167.75
$ws.Range("H99").Value = 3526.4707
$ws.Range("I99").Value = 3944.8667
$ws.Range("K99").Value = 3944.8667
$ws.Range("M99").Value = -2446.8667
$ws.Range("H134").Value = 5911.4546
$ws.Range("I134").Value = 5713.359
$ws.Range("J134").Value = 6394.3125
$ws.Range("K134").Value = 17140.077
$ws.Range("L134").Value = 19182.9375
$ws.Range("M134").Value = -14605.077
$ws.Range("N134").Value = -24252.9375

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 15124.875
$ws.Range("J4").Value = 19999.834
$ws.Range("L4").Value = 19999.834
$ws.Range("N4").Value = -20223.834
$ws.Range("H9").Value = 64611.54
$ws.Range("J9").Value = 64611.54
$ws.Range("L9").Value = 64611.54
$ws.Range("N9").Value = -64947.54
$ws.Range("H10").Value = 820.2
$ws.Range("I10").Value = 400.25
$ws.Range("J10").Value = 2500
$ws.Range("K10").Value = 400.25
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = -261.25
$ws.Range("N10").Value = -2778
$ws.Range("H13").Value = 5992.6665
$ws.Range("J13").Value = 5992.6665
$ws.Range("L13").Value = 5992.6665
$ws.Range("N13").Value = -6270.6665
$ws.Range("H31").Value = 2856.7124
$ws.Range("I31").Value = 1990.2979
$ws.Range("K31").Value = 1990.2979
$ws.Range("M31").Value = -1695.2979
$ws.Range("H34").Value = 2856.7124
$ws.Range("I34").Value = 1990.2979
$ws.Range("K34").Value = 1990.2979
$ws.Range("M34").Value = -1788.2979
$ws.Range("H41").Value = 20354.5
$ws.Range("J41").Value = 38959
$ws.Range("L41").Value = 38959
$ws.Range("N41").Value = -39815
$ws.Range("H58").Value = 2127.6843
$ws.Range("I58").Value = 1722.1333
$ws.Range("K58").Value = 1722.1333
$ws.Range("M58").Value = -1519.1333
$ws.Range("H60").Value = 26500
$ws.Range("J60").Value = 45000
$ws.Range("L60").Value = 45000
$ws.Range("N60").Value = -46022
$ws.Range("H62").Value = 7443.3335
$ws.Range("I62").Value = 6981
$ws.Range("K62").Value = 6981
$ws.Range("M62").Value = -6357
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 7443.3335
$ws.Range("I65").Value = 6981
$ws.Range("K65").Value = 34905
$ws.Range("M65").Value = -31785
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H99").Value = 1964.3
$ws.Range("I99").Value = 1899.2667
$ws.Range("K99").Value = 1899.2667
$ws.Range("M99").Value = -401.2666999999999
$ws.Range("H126").Value = 1964.3
$ws.Range("I126").Value = 1899.2667
$ws.Range("K126").Value = 5697.800099999999
$ws.Range("M126").Value = -3227.800099999999
$ws.Range("H132").Value = 3408
$ws.Range("I132").Value = 3123.8
$ws.Range("J132").Value = 4039.5557
$ws.Range("K132").Value = 9371.400000000001
$ws.Range("L132").Value = 12118.6671
$ws.Range("M132").Value = -6841.400000000001
$ws.Range("N132").Value = -17178.6671
$ws.Range("H134").Value = 3992.7297
$ws.Range("I134").Value = 3908.5
$ws.Range("J134").Value = 4353.7144
$ws.Range("K134").Value = 11725.5
$ws.Range("L134").Value = 13061.1432
$ws.Range("M134").Value = -9190.5
$ws.Range("N134").Value = -18131.1432
$ws.Range("H136").Value = 2127.6843
$ws.Range("I136").Value = 1722.1333
$ws.Range("K136").Value = 5166.3999
$ws.Range("M136").Value = -2616.3999

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 40877900
$ws.Range("I4").Value = 51046776
$ws.Range("K4").Value = 153140328
$ws.Range("M4").Value = -153140216
$ws.Range("H38").Value = 244.95653
$ws.Range("J38").Value = 359.22223
$ws.Range("L38").Value = 1077.66669
$ws.Range("N38").Value = -1771.66669
$ws.Range("H80").Value = 2998.25
$ws.Range("I80").Value = 2997
$ws.Range("J80").Value = 2998.6667
$ws.Range("K80").Value = 8991
$ws.Range("L80").Value = 8996.000100000001
$ws.Range("M80").Value = -8055
$ws.Range("N80").Value = -10868.0001
$ws.Range("H83").Value = 2998.25
$ws.Range("I83").Value = 2997
$ws.Range("J83").Value = 2998.6667
$ws.Range("K83").Value = 26973
$ws.Range("L83").Value = 26988.0003
$ws.Range("M83").Value = -22293
$ws.Range("N83").Value = -36348.0003
$ws.Range("H107").Value = 895.1579
$ws.Range("I107").Value = 598
$ws.Range("J107").Value = 974.4
$ws.Range("K107").Value = 1794
$ws.Range("L107").Value = 2923.2
$ws.Range("M107").Value = 126
$ws.Range("N107").Value = -6763.2
$ws.Range("H113").Value = 2010.6471
$ws.Range("J113").Value = 2152.5386
$ws.Range("L113").Value = 6457.6158
$ws.Range("N113").Value = -10797.6158
$ws.Range("H131").Value = 3926.761
$ws.Range("J131").Value = 3868.1628
$ws.Range("L131").Value = 11604.4884
$ws.Range("N131").Value = -21684.4884
$ws.Range("H140").Value = 4471067
$ws.Range("I140").Value = 11366376
$ws.Range("J140").Value = 9396.235
$ws.Range("K140").Value = 34099128
$ws.Range("L140").Value = 28188.705
$ws.Range("M140").Value = -34093948
$ws.Range("N140").Value = -38548.705

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 6124.75
$ws.Range("I22").Value = 2749.5
$ws.Range("J22").Value = 9500
$ws.Range("K22").Value = 2749.5
$ws.Range("L22").Value = 9500
$ws.Range("M22").Value = -2220.5
$ws.Range("N22").Value = -10558
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H80").Value = 2425
$ws.Range("I80").Value = 2425
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2425
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1427
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 2425
$ws.Range("I83").Value = 2425
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12125
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -7133
$ws.Range("N83").ClearContents()
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H126").Value = 9996
$ws.Range("I126").Value = 6995.8
$ws.Range("K126").Value = 20987.4
$ws.Range("M126").Value = -18517.4
$ws.Range("H132").Value = 1496.5883
$ws.Range("I132").Value = 1576.1333
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 4728.3999
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -2198.3999
$ws.Range("N132").Value = -7760
$ws.Range("H136").Value = 38919.645
$ws.Range("J136").Value = 38919.645
$ws.Range("L136").Value = 116758.935
$ws.Range("N136").Value = -121858.935

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2173.7827
$ws.Range("I40").Value = 1721.8462
$ws.Range("K40").Value = 1721.8462
$ws.Range("M40").Value = -1585.8462
$ws.Range("H61").Value = 1870.5714
$ws.Range("I61").Value = 1870.5714
$ws.Range("K61").Value = 1870.5714
$ws.Range("M61").Value = -1668.5714
$ws.Range("H68").Value = 2494.5
$ws.Range("I68").Value = 2493.5
$ws.Range("J68").Value = 2499.5
$ws.Range("K68").Value = 2493.5
$ws.Range("L68").Value = 2499.5
$ws.Range("M68").Value = -1744.5
$ws.Range("N68").Value = -3997.5
$ws.Range("H71").Value = 2494.5
$ws.Range("I71").Value = 2493.5
$ws.Range("J71").Value = 2499.5
$ws.Range("K71").Value = 12467.5
$ws.Range("L71").Value = 12497.5
$ws.Range("M71").Value = -8723.5
$ws.Range("N71").Value = -19985.5
$ws.Range("H82").Value = 1440.5
$ws.Range("I82").Value = 1245.0869
$ws.Range("J82").Value = 2339.4
$ws.Range("K82").Value = 1245.0869
$ws.Range("L82").Value = 2339.4
$ws.Range("M82").Value = -884.0869
$ws.Range("N82").Value = -3061.4
$ws.Range("H85").Value = 1440.5
$ws.Range("I85").Value = 1245.0869
$ws.Range("J85").Value = 2339.4
$ws.Range("K85").Value = 1245.0869
$ws.Range("L85").Value = 2339.4
$ws.Range("M85").Value = 2.913099999999986
$ws.Range("N85").Value = -4835.4
$ws.Range("H113").Value = 1870.5714
$ws.Range("I113").Value = 1870.5714
$ws.Range("K113").Value = 1870.5714
$ws.Range("M113").Value = 299.4286
$ws.Range("H132").Value = 2118.5
$ws.Range("I132").Value = 1268.1111
$ws.Range("K132").Value = 3804.3333
$ws.Range("M132").Value = -1274.3333
$ws.Range("H136").Value = 1476.4865
$ws.Range("I136").Value = 944.37933
$ws.Range("J136").Value = 3405.375
$ws.Range("K136").Value = 2833.13799
$ws.Range("L136").Value = 10216.125
$ws.Range("M136").Value = -283.1379900000002
$ws.Range("N136").Value = -15316.125

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 102127
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 112239.7
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 112239.7
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -112463.7
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H62").Value = 8666.333
$ws.Range("I62").Value = 7999.5
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 7999.5
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -7375.5
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 8666.333
$ws.Range("I65").Value = 7999.5
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 39997.5
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -36877.5
$ws.Range("N65").Value = -56240
$ws.Range("H96").Value = 4421.778
$ws.Range("J96").Value = 4421.778
$ws.Range("L96").Value = 4421.778
$ws.Range("N96").Value = -7167.778
$ws.Range("H126").Value = 2837.375
$ws.Range("I126").Value = 2366.6667
$ws.Range("K126").Value = 7100.000100000001
$ws.Range("M126").Value = -4630.000100000001
$ws.Range("H132").Value = 2191.4
$ws.Range("I132").Value = 2069.8975
$ws.Range("J132").Value = 2981.1667
$ws.Range("K132").Value = 6209.6925
$ws.Range("L132").Value = 8943.500100000001
$ws.Range("M132").Value = -3679.6925
$ws.Range("N132").Value = -14003.5001
$ws.Range("H136").Value = 3674.5789
$ws.Range("I136").Value = 3525.75
$ws.Range("K136").Value = 10577.25
$ws.Range("M136").Value = -8027.25
